$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.041.43"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.086.18"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.70%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.92"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.21"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.65%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.087.66"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.60%  "

$ws.Range("E9").Value = "  -0.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.139"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.28"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.441"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000235"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.70"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.601.07"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.119"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.378.93"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.089.84"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.40"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "448.74"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.58"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.675"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.83%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.39"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.94%  "

$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.88"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.35%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.61"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.02"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.47"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.51"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.05"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0980"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.33"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.986"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.69"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "50.55"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0696"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0379"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.93"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.66%  "

$ws.Range("E41").Value = "  -1.55%  "

$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "380.46"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -7.59%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.54"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.707.90"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "125.70"
$ws.Range("D46").ClearFormats()

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.242"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.75"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.03"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.79%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.109"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.18"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.38%  "
